$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update individual job-spec text cells that were reworded.
# ---------------------------------------------------------------------------
$ws.Range("B5").Value  = "Master (Upload all team member's work due to merge conflicts)"
$ws.Range("A6").Value  = "2 Games (Memory)"
$ws.Range("C7").Value  = "Get list of users signed up to website displayed in-app using SlashDB API"

# ---------------------------------------------------------------------------
# 2. The old merged "Thesis" banner (A12:D12) is removed first so the
#    individual cells in row 12 become independently addressable again.
# ---------------------------------------------------------------------------
$ws.Range("A12:D12").UnMerge()

# ---------------------------------------------------------------------------
# 3. Row 11 gets new content; row 12 is repurposed to carry what used to
#    live in row 11 (B) plus a brand new ScreenCast task (C).
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = "Associating app with Windows store/Testing using App Cert Kit"
$ws.Range("C11").Value = "Selenium Test Suit"

$ws.Range("A12").ClearContents()
$ws.Range("B12").Value = "Get scores from UWP Application using JSON within website"
$ws.Range("C12").Value = "ScreenCast of entire project using ShareX"
$ws.Range("D12").ClearContents()

$ws.Range("C13").ClearContents()
$ws.Range("A13").Value = "Thesis"
$ws.Range("B13").ClearContents()
$ws.Range("D13").ClearContents()
$ws.Range("A13").HorizontalAlignment = -4108   # xlCenter - only A13,B13,D13 carry the centered style; C13 stays unstyled/empty
$ws.Range("B13").HorizontalAlignment = -4108
$ws.Range("D13").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 4. New thesis sub-task row 14.
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "System Design Section with Architechture of project"
$ws.Range("B14").Value = "Methodology Section explaining Agile as preferred methodology"
$ws.Range("D14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("C14").HorizontalAlignment = -4108        # xlCenter (empty styled cell)

# ---------------------------------------------------------------------------
# 5. Row 15 content changes.
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "Input on all sections referencing above jobs"
$ws.Range("B15").Value = "Input on all sections referencing above jobs"
$ws.Range("C15").Value = "System Evaluation section/Testing/Requirements"
$ws.Range("D15").Value = "Input on all sections referencing above jobs"

# ---------------------------------------------------------------------------
# 6. New rows 16-19.
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = "Write about above in Technology Review"
$ws.Range("B16").Value = "Write about above in Technology Review"
$ws.Range("C16").Value = "Input on all sections referencing above jobs"
$ws.Range("D16").Value = "Write about above in Technology Review"

$ws.Range("B17").Value = "Proof-read thesis, correct grammatical errors/unreadable sentences"
$ws.Range("C17").Value = "Write about above in Technology Review"
$ws.Range("D17").Value = "System Design Section with Architechture of project"

$ws.Range("D18").Value = "Referencing"

$ws.Range("A19").Value = "Conclusion"
$ws.Range("B19").Value = "Conclusion"
$ws.Range("C19").Value = "Conclusion"
$ws.Range("D19").Value = "Conclusion"

# ---------------------------------------------------------------------------
# 7. Column widths were widened slightly.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 46.666666666666664
$ws.Columns.Item(2).ColumnWidth = 59.5
$ws.Columns.Item(4).ColumnWidth = 52.833333333333336

# ---------------------------------------------------------------------------
# 8. View state: scroll position + selection moved.
# ---------------------------------------------------------------------------
$ws.Range("B18").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
